# Weekly refresh of the "Poroto verde" series: a new observation is
# inserted at row 392 (pushing the existing 392-434 rows down to 393-435),
# and the worksheet's used range grows from A1:R434 to A1:R435.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 392, shifting every row
# below it (392-434) down by one (to 393-435).
$ws.Rows(392).Insert()

# Populate the newly inserted row 392 with the new data point.
$ws.Cells.Item(392, 1).Value  = 8
$ws.Cells.Item(392, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(392, 3).Value  = "Coquimbo"
$ws.Cells.Item(392, 4).Value  = 45154
$ws.Cells.Item(392, 5).Value  = 4
$ws.Cells.Item(392, 6).Value  = 100112031
$ws.Cells.Item(392, 7).Value  = "Poroto verde"
$ws.Cells.Item(392, 8).Value  = "Magnum"
$ws.Cells.Item(392, 9).Value  = "Primera"
$ws.Cells.Item(392, 10).Value = 400
$ws.Cells.Item(392, 11).Value = 30000
$ws.Cells.Item(392, 12).Value = 31000
$ws.Cells.Item(392, 13).Value = 30500
$ws.Cells.Item(392, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(392, 15).Value = "Perú"
$ws.Cells.Item(392, 16).Value = 1220
$ws.Cells.Item(392, 17).Value = 25
$ws.Cells.Item(392, 18).Value = "Hortaliza"
